$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.967.73"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3
$ws.Range("D3").Value = "1.955.87"
$ws.Range("E3").Value = "  -0.57%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.37%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.03"
$ws.Range("E5").Value = "  -1.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4846"
$ws.Range("E7").Value = "  +0.28%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2938"
$ws.Range("E8").Value = "  -0.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07022"
$ws.Range("E9").Value = "  +3.29%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").Value = "  +2.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.20"
$ws.Range("E11").Value = "  -2.19%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07803"
$ws.Range("E12").Value = "  +0.76%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.941.93"
$ws.Range("E13").Value = "  -1.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.450"
$ws.Range("E14").Value = "  -0.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6988"
$ws.Range("E15").Value = "  +0.45%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "280.15"
$ws.Range("E16").Value = "  -4.03%  "

# Row 17
$ws.Range("D17").Value = "30.981.67"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.30"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007795"
$ws.Range("E19").Value = "  +0.91%  "

# Row 20
$ws.Range("D20").Value = "2.209.42"
$ws.Range("E20").Value = "  -0.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.546"
$ws.Range("E22").Value = "  -2.32%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.488"
$ws.Range("E24").Value = "  -2.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.811"
$ws.Range("E25").Value = "  -1.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.66"
$ws.Range("E26").Value = "  -1.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.78"
$ws.Range("E27").Value = "  -1.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.171"
$ws.Range("E28").Value = "  -0.51%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1047"
$ws.Range("E29").Value = "  -2.21%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.386"
$ws.Range("E30").Value = "  -4.07%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.570"
$ws.Range("E31").Value = "  -1.91%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.591"
$ws.Range("E32").Value = "  -3.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.435"
$ws.Range("E33").Value = "  -1.32%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04883"
$ws.Range("E34").Value = "  -4.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7451"
$ws.Range("E35").Value = "  -3.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.165"
$ws.Range("E36").Value = "  -0.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.737"
$ws.Range("E37").Value = "  +0.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02000"
$ws.Range("E38").Value = "  -2.54%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.686"
$ws.Range("E39").Value = "  -1.50%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.516"
$ws.Range("E40").Value = "  +0.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.48"
$ws.Range("E41").Value = "  +8.19%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.120"
$ws.Range("E42").Value = "  -0.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9016"
$ws.Range("E43").Value = "  +1.56%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.99"
$ws.Range("E44").Value = "  -1.16%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4438"
$ws.Range("E45").Value = "  -0.48%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.007"
$ws.Range("E46").Value = "  +6.59%  "

# Row 47
$ws.Range("E47").Value = "  +0.22%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "990.42"
$ws.Range("E48").Value = "  +6.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.304"
$ws.Range("E49").Value = "  -1.09%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1245"
$ws.Range("E50").Value = "  -2.46%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.91"
$ws.Range("E51").Value = "  -0.28%  "
